# Generate Report for Handoff
# b.md has now been handed back (zh-cn and de-de) and is ready for handoff.
# Update the Overview sheet and the per-locale (zh-cn / de-de) sheets for the
# b.md row (row 3) to reflect the new status, datetimes, handback files and
# the "content duplicate" / error-detail columns, mirroring how a.md's row
# already looks after its own handback cycle.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 == b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-27 06:34:59"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 == b.md
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-27 06:34:55"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d044e13b9ed55763fc0875923c122f5d1763a7a4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3bf904b4c25f9b254c5856100c5f1e390a843065/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------
# de-de sheet: row 3 == b.md
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-27 06:34:59"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d044e13b9ed55763fc0875923c122f5d1763a7a4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3bf904b4c25f9b254c5856100c5f1e390a843065/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.14
